# Insert a new weekly data row at row 346 (pushes existing rows 346-358 down to 347-359)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(346).Insert()

$ws.Cells.Item(346, 1).Value = 8
$ws.Cells.Item(346, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(346, 3).Value = 'Coquimbo'
$ws.Cells.Item(346, 4).Value = 44939
$ws.Cells.Item(346, 5).Value = 4
$ws.Cells.Item(346, 6).Value = 100112012
$ws.Cells.Item(346, 7).Value = 'Espinaca'
$ws.Cells.Item(346, 8).Value = 'Sin especificar'
$ws.Cells.Item(346, 9).Value = 'Primera'
$ws.Cells.Item(346, 10).Value = 1600
$ws.Cells.Item(346, 11).Value = 500
$ws.Cells.Item(346, 12).Value = 600
$ws.Cells.Item(346, 13).Value = 550
$ws.Cells.Item(346, 14).Value = '$/atado 300 a 500 gramos'
$ws.Cells.Item(346, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(346, 16).Value = 1100
$ws.Cells.Item(346, 17).Value = 0.5
$ws.Cells.Item(346, 18).Value = 'Hortaliza'
